# feat: add 2022-Q4 data
#
# Original workbook has two sheets: "总计" (summary) and "2022-Q3" (fund
# holdings for Q3). This change:
#   1. Keeps the existing "2022-Q3" sheet's data+formatting exactly as-is,
#      but it becomes the NEW sheet named "2022-Q4" with a few of its
#      numbers updated (the data provider re-used the same sheet for the
#      new quarter then cloned off the old one).
#   2. Adds a brand-new sheet (after the Q4 sheet) that is a fresh copy
#      holding the original, unmodified "2022-Q3" figures.
#   3. Appends a "2022-Q3" row to the "总计" sheet (it used to have only one
#      data row) and rewrites the existing data row to describe Q4 instead.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)        # 总计
$wsQ3Orig = $wb.Worksheets.Item(2)   # currently "2022-Q3"

# ---------------------------------------------------------------------
# Step 1: duplicate the existing "2022-Q3" sheet so the ORIGINAL figures
# survive unchanged on their own tab, positioned right after it.
# ---------------------------------------------------------------------
$wsQ3Orig.Copy($null, $wsQ3Orig)
$wsQ3New = $wb.Worksheets.Item(3)

# Rename: the original sheet becomes "2022-Q4" (so it keeps sheetId/rId
# of the original "2022-Q3" sheet), the freshly made copy becomes the new
# "2022-Q3" (holding the untouched old data).
$wsQ3Orig.Name = "2022-Q4"
$wsQ3New.Name = "2022-Q3"
$wsQ4 = $wsQ3Orig

# ---------------------------------------------------------------------
# Step 2: update the Q4 sheet's data row (fund size / position / etc.)
# Several columns hold numbers formatted as text in the source file, so
# route the new values through a scratch cell forced to Text and paste
# only the *value* across - this keeps the same inline-string cell type
# without dragging a stray number-format style onto the target cell.
# ---------------------------------------------------------------------
function Set-TextValue($ws, $addr, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)   # xlPasteValues
    $scratch.Clear()
}

Set-TextValue $wsQ4 "D2" "15.29"
Set-TextValue $wsQ4 "E2" "76.11"
Set-TextValue $wsQ4 "F2" "0.91"
Set-TextValue $wsQ4 "G2" "0.1391"
$wsQ4.Range("H2").Value = 10

# The Q4 sheet's header row + A2 switch to the style already used by the
# "总计" sheet's header/label cells (bold, borderless-color thin border)
# instead of the original sheet's style.
$ws1.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("A2").Copy()
$wsQ4.Range("A2").PasteSpecial(-4122)

# Match the "总计" sheet's page margins (0.75in/1in/0.5in) on the Q4 tab.
$wsQ4.PageSetup.LeftMargin = 54
$wsQ4.PageSetup.RightMargin = 54
$wsQ4.PageSetup.TopMargin = 72
$wsQ4.PageSetup.BottomMargin = 72
$wsQ4.PageSetup.HeaderMargin = 36
$wsQ4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# Step 3: "总计" sheet - existing row 2 now describes 2022-Q4, and a new
# row 3 is appended with the original 2022-Q3 totals.
# ---------------------------------------------------------------------
$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("D2").Value = 0.14

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q3"
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 0.13

$ws1.Range("A2").Copy()
$ws1.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

# Keep the workbook-level active tab on "总计" (index 0), matching the
# untouched <bookViews> element from the source file.
$ws1.Activate()

